$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Fill the "No." column (A) with sequential row numbers 1..85 for rows 2..86 ---
# Copy formatting from B2, which already carries the target cell style, into A2:A86
# so the new cells end up with the same style index used throughout the sheet.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A86").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

for ($i = 2; $i -le 86; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# --- Update the sheet view: clear the old scrolled/selected state and select A2:A86 ---
$ws.Range("A2:A86").Select() | Out-Null

# --- Update the workbook window position/size recorded in the bookViews ---
$win = $excel.ActiveWindow
$win.Left = 29310
$win.Top = 675
$win.Width = 25845
$win.Height = 14850
